$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G to fit the new, longer tag strings
$ws.Columns.Item(7).ColumnWidth = 97.42578125

# Update the "Tags" column (G) with refined/more specific tag strings
# for rows whose tag text actually changes (katakana sub-series tags and
# the "special foreign" katakana rows that now carry detailed sub-tags).
$ws.Cells.Item(120, 7).Value = 'katakana gojuon seion s_series sh_series'
$ws.Cells.Item(125, 7).Value = 'katakana gojuon seion t_series ch_series'
$ws.Cells.Item(126, 7).Value = 'katakana gojuon seion t_series ts_series'
$ws.Cells.Item(136, 7).Value = 'katakana gojuon seion h_series f_series'
$ws.Cells.Item(161, 7).Value = 'katakana dakuon s_series z_series j_series'
$ws.Cells.Item(216, 7).Value = 'katakana special foreign gojuon_sp seion_sp a_series_sp y_series_x'
$ws.Cells.Item(217, 7).Value = 'katakana special foreign gojuon_sp seion_sp a_series_sp y_series_x'
$ws.Cells.Item(218, 7).Value = 'katakana special foreign gojuon_sp seion_sp a_series_sp w_series_x'
$ws.Cells.Item(219, 7).Value = 'katakana special foreign gojuon_sp seion_sp a_series_sp w_series_x'
$ws.Cells.Item(220, 7).Value = 'katakana special foreign gojuon_sp seion_sp a_series_sp w_series_x'
$ws.Cells.Item(221, 7).Value = 'katakana special foreign dakuon_sp a_series_sp v_series v_series_x'
$ws.Cells.Item(222, 7).Value = 'katakana special foreign dakuon_sp a_series_sp v_series v_series_x'
$ws.Cells.Item(223, 7).Value = 'katakana special foreign dakuon_sp a_series_sp v_series v_series_x'
$ws.Cells.Item(224, 7).Value = 'katakana special foreign dakuon_sp a_series_sp v_series v_series_x'
$ws.Cells.Item(225, 7).Value = 'katakana special foreign dakuon_sp a_series_sp v_series v_series_x'
$ws.Cells.Item(226, 7).Value = 'katakana special foreign gojuon_sp seion_sp s_series_sp s_series_x sh_series_sp sh_series_x'
$ws.Cells.Item(227, 7).Value = 'katakana special foreign dakuon_sp s_series_sp s_series_x sh_series_sp sh_series_x j_series_sp j_series_x'
$ws.Cells.Item(228, 7).Value = 'katakana special foreign gojuon_sp seion_sp t_series_sp t_series_x ch_series_sp ch_series_x'
$ws.Cells.Item(229, 7).Value = 'katakana special foreign gojuon_sp seion_sp t_series_sp t_series_x'
$ws.Cells.Item(230, 7).Value = 'katakana special foreign gojuon_sp seion_sp t_series_sp t_series_x'
$ws.Cells.Item(231, 7).Value = 'katakana special foreign dakuon_sp t_series_sp t_series_x d_series_sp d_series_x'
$ws.Cells.Item(232, 7).Value = 'katakana special foreign dakuon_sp t_series_sp t_series_x d_series_sp d_series_x'
$ws.Cells.Item(233, 7).Value = 'katakana special foreign gojuon_sp seion_sp t_series_sp t_series_x ts_series_sp ts_series_x'
$ws.Cells.Item(234, 7).Value = 'katakana special foreign gojuon_sp seion_sp t_series_sp t_series_x ts_series_sp ts_series_x'
$ws.Cells.Item(235, 7).Value = 'katakana special foreign gojuon_sp seion_sp t_series_sp t_series_x ts_series_sp ts_series_x'
$ws.Cells.Item(236, 7).Value = 'katakana special foreign gojuon_sp seion_sp t_series_sp t_series_x ts_series_sp ts_series_x'
$ws.Cells.Item(237, 7).Value = 'katakana special foreign gojuon_sp seion_sp h_series_sp f_series_sp f_series_x'
$ws.Cells.Item(238, 7).Value = 'katakana special foreign gojuon_sp seion_sp h_series_sp f_series_sp f_series_x'
$ws.Cells.Item(239, 7).Value = 'katakana special foreign gojuon_sp seion_sp h_series_sp f_series_sp f_series_x'
$ws.Cells.Item(240, 7).Value = 'katakana special foreign gojuon_sp seion_sp h_series_sp f_series_sp f_series_x'
$ws.Cells.Item(241, 7).Value = 'katakana special foreign yoon_dakuon_sp a_series_sp v_series v_series_x vy_series vy_series_x'
$ws.Cells.Item(242, 7).Value = 'katakana special foreign yoon_dakuon_sp a_series_sp v_series v_series_x vy_series vy_series_x'
$ws.Cells.Item(243, 7).Value = 'katakana special foreign yoon_dakuon_sp a_series_sp v_series v_series_x vy_series vy_series_x'
$ws.Cells.Item(244, 7).Value = 'katakana special foreign yoon_sp t_series_sp ty_series ty_series_x'
$ws.Cells.Item(245, 7).Value = 'katakana special foreign yoon_dakuon_sp t_series_sp ty_series ty_series_x'
$ws.Cells.Item(246, 7).Value = 'katakana special foreign yoon_dakuon_sp f_series_sp f_series_x fy_series fy_series_x'
$ws.Cells.Item(247, 7).Value = 'katakana special foreign gojuon_sp seion_sp s_series_sp sw_series sw_series_x'
$ws.Cells.Item(248, 7).Value = 'katakana special foreign gojuon_sp seion_sp s_series_sp s_series_x sw_series sw_series_x'
$ws.Cells.Item(249, 7).Value = 'katakana special foreign gojuon_sp seion_sp s_series_sp sw_series sw_series_x'
$ws.Cells.Item(250, 7).Value = 'katakana special foreign gojuon_sp seion_sp s_series_sp sw_series sw_series_x'
$ws.Cells.Item(251, 7).Value = 'katakana special foreign gojuon_sp seion_sp s_series_sp sw_series sw_series_x'
$ws.Cells.Item(252, 7).Value = 'katakana special foreign dakuon_sp s_series_sp zw_series zw_series_x'
$ws.Cells.Item(253, 7).Value = 'katakana special foreign dakuon_sp s_series_sp z_series_x zw_series zw_series_x'
$ws.Cells.Item(254, 7).Value = 'katakana special foreign dakuon_sp s_series_sp zw_series zw_series_x'
$ws.Cells.Item(255, 7).Value = 'katakana special foreign dakuon_sp s_series_sp zw_series zw_series_x'
$ws.Cells.Item(256, 7).Value = 'katakana special foreign dakuon_sp s_series_sp zw_series zw_series_x'
